$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$hdr = $ws.Range("A1:U1")
$scratch = $ws.Range("A100:U100")
$hdr.Copy()
$scratch.PasteSpecial(-4122)
$hdr.ClearFormats()
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U55"), $null, 1)
$tbl.TableStyle = ""
$scratch.Copy()
$hdr.PasteSpecial(-4122)
$scratch.ClearFormats()
$scratch.ClearContents()
Write-Output "done"
